# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value mapping for sheet "展览" (rows 2-13)
$exhibitionUpdates = @{
    2  = 70
    3  = 801
    5  = 61
    6  = 103
    7  = 337
    8  = 4124
    9  = 91
    10 = 4846
    11 = 538
    12 = 1216
    13 = 83
}

# Row -> new value mapping for sheet "全部类型" (rows 2-14)
$allTypesUpdates = @{
    2  = 70
    3  = 801
    5  = 61
    6  = 103
    8  = 337
    9  = 4124
    10 = 91
    11 = 4846
    12 = 538
    13 = 1216
    14 = 83
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
